$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 87, shifting existing rows 87-101 down to 88-102.
$ws.Rows.Item(87).Insert()

# Row 87 becomes a copy of the original row 86 data (before this edit's changes).
$ws.Range("A87").Value = 4
$ws.Range("B87").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C87").Value = "Los Lagos"
$ws.Range("D87").Value = 44726
$ws.Range("E87").Value = 10
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100104
$ws.Range("H87").Value = "Frutos de pepita"
$ws.Range("I87").Value = 100104003
$ws.Range("J87").Value = "Membrillo"
$ws.Range("K87").Value = "Champion"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 300
$ws.Range("N87").Value = 13000
$ws.Range("O87").Value = 14000
$ws.Range("P87").Value = 13500
$ws.Range("Q87").Value = "`$/caja 18 kilos granel"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 750
$ws.Range("T87").Value = 18

# Row 86 itself gets a new date and volume.
$ws.Range("D86").Value = 44785
$ws.Range("M86").Value = 400
